# Update Leve price/profit figures across sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3027.111
$ws.Range("I111").Value = 2299.6667
$ws.Range("J111").Value = 3390.8333
$ws.Range("K111").Value = 6899.000100000001
$ws.Range("L111").Value = 10172.4999
$ws.Range("M111").Value = -3832.000100000001
$ws.Range("N111").Value = -16306.4999

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3674.6667
$ws.Range("I132").Value = 3607.0588
$ws.Range("K132").Value = 10821.1764
$ws.Range("M132").Value = -8291.1764

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1978.1
$ws.Range("I138").Value = 1978.1
$ws.Range("K138").Value = 5934.299999999999
$ws.Range("M138").Value = -794.2999999999993

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1155.6
$ws.Range("I2").Value = 960
$ws.Range("K2").Value = 960
$ws.Range("M2").Value = -847

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1241.8125
$ws.Range("I32").Value = 1243.1613
$ws.Range("K32").Value = 1243.1613
$ws.Range("M32").Value = -956.1613

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1155.6
$ws.Range("I116").Value = 960
$ws.Range("K116").Value = 960
$ws.Range("M116").Value = 1334

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2967.2432
$ws.Range("I132").Value = 3011.611
$ws.Range("K132").Value = 9034.832999999999
$ws.Range("M132").Value = -6504.832999999999

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1155.6
$ws.Range("I3").Value = 960
$ws.Range("K3").Value = 960
$ws.Range("M3").Value = -846

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3178.6667
$ws.Range("I20").Value = 3213.625
$ws.Range("J20").Value = 2899
$ws.Range("K20").Value = 3213.625
$ws.Range("L20").Value = 2899
$ws.Range("M20").Value = -2966.625
$ws.Range("N20").Value = -3393

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 179.5
$ws.Range("I22").Value = 192.4
$ws.Range("J22").Value = 147.25
$ws.Range("K22").Value = 192.4
$ws.Range("L22").Value = 147.25
$ws.Range("M22").Value = -19.40000000000001
$ws.Range("N22").Value = -493.25

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 59083.75
$ws.Range("J81").Value = 59083.75
$ws.Range("L81").Value = 59083.75
$ws.Range("N81").Value = -61205.75

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 59083.75
$ws.Range("J84").Value = 59083.75
$ws.Range("L84").Value = 177251.25
$ws.Range("N84").Value = -187859.25

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1319.3334
$ws.Range("I16").Value = 1246.75
$ws.Range("K16").Value = 1246.75
$ws.Range("M16").Value = -959.75

# CRP row 18
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4253.7
$ws.Range("I36").Value = 2512.3333
$ws.Range("K36").Value = 2512.3333
$ws.Range("M36").Value = -2124.3333

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 4253.7
$ws.Range("I40").Value = 2512.3333
$ws.Range("K40").Value = 2512.3333
$ws.Range("M40").Value = -2352.3333

# CRP row 45
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 4900
$ws.Range("I45").Value = 4900
$ws.Range("K45").Value = 4900
$ws.Range("M45").Value = -4307

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 484.75
$ws.Range("I105").Value = 484.75
$ws.Range("K105").Value = 484.75
$ws.Range("M105").Value = 1262.25

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1319.3334
$ws.Range("I113").Value = 1246.75
$ws.Range("K113").Value = 1246.75
$ws.Range("M113").Value = 923.25

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3218
$ws.Range("I132").Value = 3222.5
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 9667.5
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -7137.5
$ws.Range("N132").Value = -14660

# CUL row 10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 57
$ws.Range("I10").Value = 57
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 171
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -32
$ws.Range("N10").ClearContents()

# CUL row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 21
$ws.Range("I13").Value = 22.2
$ws.Range("K13").Value = 66.59999999999999
$ws.Range("M13").Value = 101.4

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1095
$ws.Range("J26").Value = 395.5
$ws.Range("L26").Value = 1186.5
$ws.Range("N26").Value = -1762.5

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1899.5
$ws.Range("I87").Value = 1899.5
$ws.Range("K87").Value = 5698.5
$ws.Range("M87").Value = -4450.5

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 1899.5
$ws.Range("I90").Value = 1899.5
$ws.Range("K90").Value = 17095.5
$ws.Range("M90").Value = -10855.5

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 603.4
$ws.Range("I113").Value = 540
$ws.Range("J113").Value = 610.44446
$ws.Range("K113").Value = 1620
$ws.Range("L113").Value = 1831.33338
$ws.Range("M113").Value = 550
$ws.Range("N113").Value = -6171.33338

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2359.6924
$ws.Range("I131").Value = 2828
$ws.Range("J131").Value = 1813.3334
$ws.Range("K131").Value = 8484
$ws.Range("L131").Value = 5440.0002
$ws.Range("M131").Value = -3444
$ws.Range("N131").Value = -15520.0002

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2308.5
$ws.Range("I132").Value = 2345.7896
$ws.Range("K132").Value = 7037.3688
$ws.Range("M132").Value = -4507.3688

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 785.2857
$ws.Range("I16").Value = 807.8333
$ws.Range("J16").Value = 650
$ws.Range("K16").Value = 807.8333
$ws.Range("L16").Value = 650
$ws.Range("M16").Value = -637.8333
$ws.Range("N16").Value = -990

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2349
$ws.Range("J22").Value = 2199
$ws.Range("L22").Value = 2199
$ws.Range("N22").Value = -2789

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2349
$ws.Range("J27").Value = 2199
$ws.Range("L27").Value = 2199
$ws.Range("N27").Value = -2413

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3299.4
$ws.Range("I68").Value = 2831.3333
$ws.Range("J68").Value = 4001.5
$ws.Range("K68").Value = 2831.3333
$ws.Range("L68").Value = 4001.5
$ws.Range("M68").Value = -2082.3333
$ws.Range("N68").Value = -5499.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3299.4
$ws.Range("I71").Value = 2831.3333
$ws.Range("J71").Value = 4001.5
$ws.Range("K71").Value = 14156.6665
$ws.Range("L71").Value = 20007.5
$ws.Range("M71").Value = -10412.6665
$ws.Range("N71").Value = -27495.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 41668644
$ws.Range("I136").Value = 1747.5714
$ws.Range("J136").Value = 100002296
$ws.Range("K136").Value = 5242.7142
$ws.Range("L136").Value = 300006888
$ws.Range("M136").Value = -2692.7142
$ws.Range("N136").Value = -300011988

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1111.8572
$ws.Range("I100").Value = 1111.8572
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2223.7144
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1682.7144
$ws.Range("N100").ClearContents()

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2882.5881
$ws.Range("I122").Value = 2882.5881
$ws.Range("K122").Value = 8647.764299999999
$ws.Range("M122").Value = -6197.764299999999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6484.75
$ws.Range("I132").Value = 4812.1665
$ws.Range("K132").Value = 14436.4995
$ws.Range("M132").Value = -11906.4995
